$d = $word.ActiveDocument

# Replace the placeholder ID text and remove the trailing space run that
# followed it (Find/Replace across both runs merges them into one run).
$d.Content.Find.Execute("**ID__AFFARS_5332_topic_9__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5332_501_2__ID**", 2)

# Give the first paragraph a (borderless) paragraph border with 5-twip
# spacing on all sides, and bump its left indent from 120 to 225 twips
# (11.25pt).
$p = $d.Paragraphs(1)
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromRight = 5
$p.Format.LeftIndent = 11.25
